$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the core data row (row 2) with new values for the TimelyCare module
$ws.Range("A2").Value = "PMH_TimelyCare"
$ws.Range("B2").Value = "TimelyCare"
$ws.Range("C2").Value = "PMHdata()"
$ws.Range("D2").Value = "TimelyCare"
$ws.Range("E2").Value = "F"
$ws.Range("F2").Value = "T"
$ws.Range("G2").Value = "F"

# Add literal copies of the generated formula strings below (rows 10-12)
$ws.Range("B10").Value = "mod_Accordion_ui('PMH_TimelyCare')"
$ws.Range("B11").Value = "mod_Accordion_server('PMH_TimelyCare', selector=selection, data=PMHdata(), title = c('TimelyCare'), Visible = T)"
$ws.Range("B12").Value = "mod_info_server('PMH_TimelyCare', selector = selection, data = PMHdata(), rownametitle = c('TimelyCare'), phone = F, website = T, email = F)"

# Match the recorded selection change
$ws.Range("B10").Select() | Out-Null
